$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C1").Value = "nickname"
$ws.Range("C2").Value = "aa's nickname"
$ws.Range("C3").Value = "bb's nickname"
$ws.Range("C4").Value = "cc's nickname"
$ws.Range("C5").Value = "dd's nickname"
